$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'70.525.55"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -2.08%  "

$ws.Range("D3").Value = "'2.551.82"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -4.92%  "

$ws.Range("E4").Value = "  +0.06%  "

$ws.Range("D5").Value = "'577.90"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -3.38%  "

$ws.Range("D6").Value = "'169.92"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -2.56%  "

$ws.Range("E7").Value = "  +0.09%  "

$ws.Range("D8").Value = "'0.510"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -2.56%  "

$ws.Range("D9").Value = "'2.550.09"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -4.93%  "

$ws.Range("D10").Value = "'0.166"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.48%  "

$ws.Range("D11").Value = "'0.169"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.15%  "

$ws.Range("E12").Value = "  -3.03%  "

$ws.Range("E13").Value = "  -3.24%  "

$ws.Range("D14").Value = "'3.018.38"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -4.89%  "

$ws.Range("D15").Value = "'0.0000181"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.79%  "

$ws.Range("D16").Value = "'70.443.69"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -2.04%  "

$ws.Range("D17").Value = "'25.12"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -4.00%  "

$ws.Range("D18").Value = "'2.559.56"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -4.49%  "

$ws.Range("D19").Value = "'11.71"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -3.78%  "

$ws.Range("D20").Value = "'362.36"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -2.66%  "

$ws.Range("D21").Value = "'7.50"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -7.75%  "

$ws.Range("D22").Value = "'3.97"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -4.75%  "

$ws.Range("D23").Value = "'2.01"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.01%  "

$ws.Range("E24").Value = "  +0.01%  "

$ws.Range("D25").Value = "'70.03"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -3.24%  "

$ws.Range("D26").Value = "'4.09"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -5.66%  "

$ws.Range("D27").Value = "'9.32"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -4.84%  "

$ws.Range("D28").Value = "'2.690.21"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -4.66%  "

$ws.Range("D29").Value = "'1.00"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.06%  "

$ws.Range("D30").Value = "'0.0₃0927"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -4.93%  "

$ws.Range("D31").Value = "'7.88"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -1.94%  "

$ws.Range("D32").Value = "'485.90"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -3.36%  "

$ws.Range("E33").Value = "  +0.32%  "

$ws.Range("E34").Value = "  -2.88%  "

$ws.Range("E35").Value = "  +0.00%  "

$ws.Range("B36").Value = "Kaspa"
$ws.Range("C36").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D36").Value = "'0.117"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +7.28%  "

$ws.Range("B37").Value = "Monero"
$ws.Range("C37").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D37").Value = "'156.95"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -4.33%  "

$ws.Range("D38").Value = "'18.71"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -4.74%  "

$ws.Range("D39").Value = "'18.84"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.34%  "

$ws.Range("E40").Value = "  -3.21%  "

$ws.Range("E41").Value = "  +0.02%  "

$ws.Range("D42").Value = "'1.68"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -5.05%  "

$ws.Range("E43").Value = "  -4.70%  "

$ws.Range("D44").Value = "'2.47"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -3.86%  "

$ws.Range("D45").Value = "'0.321"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -3.41%  "

$ws.Range("D46").Value = "'38.45"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -2.46%  "

$ws.Range("D47").Value = "'145.12"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -7.24%  "

$ws.Range("D48").Value = "'3.56"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -4.46%  "

$ws.Range("D49").Value = "'0.532"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -5.51%  "

$ws.Range("D50").Value = "'1.64"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -6.39%  "

$ws.Range("D51").Value = "'0.595"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -2.00%  "
